$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.593.83'
$ws.Range('E2').Value = '  -3.32%  '

$ws.Range('D3').Value = '1.849.41'
$ws.Range('E3').Value = '  -3.87%  '

$ws.Range('D4').Value = '''0.9998'
$ws.Range('E4').Value = '  -1.41%  '

$ws.Range('D5').Value = '''335.98'
$ws.Range('E5').Value = '  +2.99%  '

$ws.Range('D6').Value = '''1.0000'
$ws.Range('E6').Value = '  -1.23%  '

$ws.Range('D7').Value = '''0.4668'
$ws.Range('E7').Value = '  -3.02%  '

$ws.Range('D8').Value = '''0.3907'
$ws.Range('E8').Value = '  -3.65%  '

$ws.Range('D9').Value = '''46.16'
$ws.Range('E9').Value = '  -2.70%  '

$ws.Range('D10').Value = '''0.07904'
$ws.Range('E10').Value = '  -3.59%  '

$ws.Range('D11').Value = '''0.9791'
$ws.Range('E11').Value = '  -2.78%  '

$ws.Range('D12').Value = '''22.29'
$ws.Range('E12').Value = '  -5.88%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.864.90'
$ws.Range('E13').Value = '  -1.72%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''5.819'
$ws.Range('E14').Value = '  -4.41%  '

$ws.Range('D15').Value = '''6.975'
$ws.Range('E15').Value = '  -4.26%  '

$ws.Range('D16').Value = '''0.06937'
$ws.Range('E16').Value = '  +0.87%  '

$ws.Range('D17').Value = '''87.67'
$ws.Range('E17').Value = '  -4.19%  '

$ws.Range('D18').Value = '''1.000'
$ws.Range('E18').Value = '  -1.38%  '

$ws.Range('D19').Value = '''0.000009993'
$ws.Range('E19').Value = '  -3.77%  '

$ws.Range('E20').Value = '  -3.23%  '

$ws.Range('E21').Value = '  -1.13%  '

$ws.Range('D22').Value = '28.616.64'
$ws.Range('E22').Value = '  -3.23%  '

$ws.Range('D23').Value = '''5.394'
$ws.Range('E23').Value = '  -4.79%  '

$ws.Range('D24').Value = '''11.26'
$ws.Range('E24').Value = '  -5.74%  '

$ws.Range('D25').Value = '''2.150'
$ws.Range('E25').Value = '  -1.52%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''152.71'
$ws.Range('E26').Value = '  -2.16%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''19.43'
$ws.Range('E27').Value = '  -2.93%  '

$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''6.067'
$ws.Range('E28').Value = '  -4.74%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '''2.019'
$ws.Range('E29').Value = '  -3.28%  '

$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''117.41'
$ws.Range('E30').Value = '  -2.70%  '

$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''0.9720'
$ws.Range('E31').Value = '  -3.60%  '

$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.09351'
$ws.Range('E32').Value = '  -2.56%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '''5.364'
$ws.Range('E33').Value = '  -4.15%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.471'
$ws.Range('E34').Value = '  -2.67%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '''1.349'
$ws.Range('E35').Value = '  -2.57%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.06141'
$ws.Range('E36').Value = '  -5.69%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02201'
$ws.Range('E37').Value = '  -3.42%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.170'
$ws.Range('E38').Value = '  -2.27%  '

$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''7.690'
$ws.Range('E39').Value = '  -2.11%  '

$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5708'
$ws.Range('E40').Value = '  -3.75%  '

$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').Value = '''10.13'
$ws.Range('E41').Value = '  -5.76%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '''2.444'
$ws.Range('E42').Value = '  -2.72%  '

$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').Value = '''0.1794'
$ws.Range('E43').Value = '  -2.58%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''1.243'
$ws.Range('E44').Value = '  -3.03%  '

$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''11.80'
$ws.Range('E45').Value = '  -4.89%  '

$ws.Range('D46').Value = '''0.5364'
$ws.Range('E46').Value = '  -3.19%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.07099'
$ws.Range('E47').Value = '  -5.33%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.908'
$ws.Range('E48').Value = '  -2.49%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''112.96'
$ws.Range('E49').Value = '  -4.57%  '

$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').Value = '''2.355'
$ws.Range('E50').Value = '  -3.35%  '

$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').Value = '''0.9990'
$ws.Range('E51').Value = '  -1.36%  '
